$d = $word.ActiveDocument

# 1. Update the "Hello!" paragraph text to mention switching into software development.
$d.Content.Find.Execute(
    "living in London. I appreciate",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "living in London looking to make the switch into software development. I appreciate",
    2
)

# 2. Remove one of the three consecutive empty paragraphs that sit between the
#    "Before starting my journey..." paragraph and the "EXPERIENCE" heading.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "Before starting my journey with tech and technical recruitment in 2017 I worked for a few manufacturing and warehousing companies from February 2011 until May 2017 primarily focusing on quality control, inventory auditing and inventory processing.") {
        $emptyPara = $d.Paragraphs.Item($i + 1)
        $emptyPara.Range.Delete()
        break
    }
}
